# Slide 4, "Content Placeholder 2" shape: first paragraph currently reads
#   "Create a method that checks if the integer value is an even or odd
#    number, return 1 for even and 0 for odd. "
# The commit swaps the two numbers (1 <-> 0) and, in doing so, splits the
# single run into five runs (matching how PowerPoint itself fragments a
# run when only parts of it are retyped).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

# First paragraph of the text box.
$para1 = $tr.Paragraphs(1, 1)

# Locate the two single-character tokens that need to swap, scoped to
# this paragraph only (1-based, paragraph-relative offsets).
$oldText = $para1.Text
$iReturn1 = $oldText.IndexOf("return 1")
$iAnd0    = $oldText.IndexOf("and 0")

# "return 1 " -> "return 0 " (9 chars, starts right after "Create a
# method ... number, ").
$segReturn = $para1.Characters($iReturn1 + 1, 9)
$segReturn.Text = "return 0 "

# "and 0 " -> "and 1 " (6 chars).
$segAnd = $para1.Characters($iAnd0 + 1, 6)
$segAnd.Text = "and 1 "
